$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

$ws.Range("A100").Value = "Assigned_User_Button"
$ws.Range("B100").Value = "//a[normalize-space()='Assigned Users']"
$ws.Range("C100").Value = "By.xpath"

$ws.Range("B101").Value = "//button[normalize-space()='Collect Garbage']"
$ws.Range("A101").Value = "Collect_Garbage_Button"
$ws.Range("C101").Value = "By.xpath"

$ws.Range("B102").Value = "//select[@class='form-select']"
$ws.Range("A102").Value = "Waste_Type_Selector"
$ws.Range("C102").Value = "By.xpath"

$ws.Range("A102").Select()
